$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1838.8313
$ws.Range("I15").Value = 1838.8313
$ws.Range("K15").Value = 5516.4939
$ws.Range("M15").Value = -5347.4939
$ws.Range("H33").Value = 118.333336
$ws.Range("I33").Value = 121.23077
$ws.Range("J33").Value = 99.5
$ws.Range("K33").Value = 121.23077
$ws.Range("L33").Value = 99.5
$ws.Range("M33").Value = 107.76923
$ws.Range("N33").Value = -557.5
$ws.Range("H98").Value = 882.86957
$ws.Range("I98").Value = 891.1818
$ws.Range("J98").Value = 700
$ws.Range("K98").Value = 891.1818
$ws.Range("L98").Value = 700
$ws.Range("M98").Value = 606.8182
$ws.Range("N98").Value = -3696
$ws.Range("H113").Value = 2682.318
$ws.Range("I113").Value = 2194.6428
$ws.Range("J113").Value = 3535.75
$ws.Range("K113").Value = 2194.6428
$ws.Range("L113").Value = 3535.75
$ws.Range("M113").Value = 1059.3572
$ws.Range("N113").Value = -10043.75
$ws.Range("H116").Value = 3249.1667
$ws.Range("I116").Value = 3501.6667
$ws.Range("J116").Value = 2996.6667
$ws.Range("K116").Value = 3501.6667
$ws.Range("L116").Value = 2996.6667
$ws.Range("M116").Value = -59.66670000000022
$ws.Range("N116").Value = -9880.6667
$ws.Range("H122").Value = 882.86957
$ws.Range("I122").Value = 891.1818
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 2673.5454
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -223.5454
$ws.Range("N122").Value = -7000
$ws.Range("H127").Value = 1283
$ws.Range("I127").Value = 1049
$ws.Range("J127").Value = 1517
$ws.Range("K127").Value = 3147
$ws.Range("L127").Value = 4551
$ws.Range("M127").Value = 1813
$ws.Range("N127").Value = -14471
$ws.Range("H130").Value = 48000
$ws.Range("J130").Value = 48000
$ws.Range("L130").Value = 48000
$ws.Range("N130").Value = -58040
$ws.Range("H132").Value = 2470.8438
$ws.Range("I132").Value = 2016.5883
$ws.Range("J132").Value = 4252.923
$ws.Range("K132").Value = 6049.7649
$ws.Range("L132").Value = 12758.769
$ws.Range("M132").Value = -3519.7649
$ws.Range("N132").Value = -17818.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 53245.26
$ws.Range("I2").Value = 83821.75
$ws.Range("J2").Value = 828.4286
$ws.Range("K2").Value = 83821.75
$ws.Range("L2").Value = 828.4286
$ws.Range("M2").Value = -83708.75
$ws.Range("N2").Value = -1054.4286
$ws.Range("H32").Value = 4170.269
$ws.Range("I32").Value = 3242.1357
$ws.Range("K32").Value = 3242.1357
$ws.Range("M32").Value = -2955.1357
$ws.Range("H61").Value = 958.96204
$ws.Range("I61").Value = 804.5
$ws.Range("J61").Value = 1743.1538
$ws.Range("K61").Value = 804.5
$ws.Range("L61").Value = 1743.1538
$ws.Range("M61").Value = -592.5
$ws.Range("N61").Value = -2167.1538
$ws.Range("H116").Value = 53245.26
$ws.Range("I116").Value = 83821.75
$ws.Range("J116").Value = 828.4286
$ws.Range("K116").Value = 83821.75
$ws.Range("L116").Value = 828.4286
$ws.Range("M116").Value = -81527.75
$ws.Range("N116").Value = -5416.4286
$ws.Range("H132").Value = 7476.5293
$ws.Range("I132").Value = 5619.1304
$ws.Range("J132").Value = 11360.182
$ws.Range("K132").Value = 16857.3912
$ws.Range("L132").Value = 34080.546
$ws.Range("M132").Value = -14327.3912
$ws.Range("N132").Value = -39140.546
$ws.Range("H134").Value = 50105.668
$ws.Range("J134").Value = 50105.668
$ws.Range("L134").Value = 50105.668
$ws.Range("N134").Value = -60245.668
$ws.Range("H136").Value = 958.96204
$ws.Range("I136").Value = 804.5
$ws.Range("J136").Value = 1743.1538
$ws.Range("K136").Value = 2413.5
$ws.Range("L136").Value = 5229.4614
$ws.Range("M136").Value = 136.5
$ws.Range("N136").Value = -10329.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 53245.26
$ws.Range("I3").Value = 83821.75
$ws.Range("J3").Value = 828.4286
$ws.Range("K3").Value = 83821.75
$ws.Range("L3").Value = 828.4286
$ws.Range("M3").Value = -83707.75
$ws.Range("N3").Value = -1056.4286
$ws.Range("H86").Value = 2415.2173
$ws.Range("I86").Value = 2047.3334
$ws.Range("J86").Value = 2816.5454
$ws.Range("K86").Value = 2047.3334
$ws.Range("L86").Value = 2816.5454
$ws.Range("M86").Value = -924.3334
$ws.Range("N86").Value = -5062.5454
$ws.Range("H89").Value = 2415.2173
$ws.Range("I89").Value = 2047.3334
$ws.Range("J89").Value = 2816.5454
$ws.Range("K89").Value = 10236.667
$ws.Range("L89").Value = 14082.727
$ws.Range("M89").Value = -4620.666999999999
$ws.Range("N89").Value = -25314.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38233.64
$ws.Range("I31").Value = 3585.027
$ws.Range("J31").Value = 91650.25
$ws.Range("K31").Value = 3585.027
$ws.Range("L31").Value = 91650.25
$ws.Range("M31").Value = -3290.027
$ws.Range("N31").Value = -92240.25
$ws.Range("H34").Value = 38233.64
$ws.Range("I34").Value = 3585.027
$ws.Range("J34").Value = 91650.25
$ws.Range("K34").Value = 3585.027
$ws.Range("L34").Value = 91650.25
$ws.Range("M34").Value = -3383.027
$ws.Range("N34").Value = -92054.25
$ws.Range("H58").Value = 803.4423
$ws.Range("I58").Value = 728.1111
$ws.Range("J58").Value = 1287.7142
$ws.Range("K58").Value = 728.1111
$ws.Range("L58").Value = 1287.7142
$ws.Range("M58").Value = -525.1111
$ws.Range("N58").Value = -1693.7142
$ws.Range("H122").Value = 1697.8
$ws.Range("I122").Value = 1697.4117
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 5092.2351
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -2642.2351
$ws.Range("N122").Value = -10000
$ws.Range("H136").Value = 803.4423
$ws.Range("I136").Value = 728.1111
$ws.Range("J136").Value = 1287.7142
$ws.Range("K136").Value = 2184.3333
$ws.Range("L136").Value = 3863.1426
$ws.Range("M136").Value = 365.6667000000002
$ws.Range("N136").Value = -8963.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 5660
$ws.Range("I57").Value = 3980
$ws.Range("K57").Value = 11940
$ws.Range("M57").Value = -11381
$ws.Range("H115").Value = 2500
$ws.Range("I115").Value = 1500
$ws.Range("J115").Value = 3500
$ws.Range("K115").Value = 4500
$ws.Range("L115").Value = 10500
$ws.Range("M115").Value = -3325
$ws.Range("N115").Value = -12850
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 536.65717
$ws.Range("I131").Value = 459.43332
$ws.Range("J131").Value = 1000
$ws.Range("K131").Value = 1378.29996
$ws.Range("L131").Value = 3000
$ws.Range("M131").Value = 3661.70004
$ws.Range("N131").Value = -13080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4815.5
$ws.Range("I122").Value = 5584.4346
$ws.Range("J122").Value = 1278.4
$ws.Range("K122").Value = 16753.3038
$ws.Range("L122").Value = 3835.2
$ws.Range("M122").Value = -14303.3038
$ws.Range("N122").Value = -8735.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 424.125
$ws.Range("I22").Value = 339.13333
$ws.Range("J22").Value = 565.7778
$ws.Range("K22").Value = 339.13333
$ws.Range("L22").Value = 565.7778
$ws.Range("M22").Value = -44.13333
$ws.Range("N22").Value = -1155.7778
$ws.Range("H27").Value = 424.125
$ws.Range("I27").Value = 339.13333
$ws.Range("J27").Value = 565.7778
$ws.Range("K27").Value = 339.13333
$ws.Range("L27").Value = 565.7778
$ws.Range("M27").Value = -232.13333
$ws.Range("N27").Value = -779.7778
$ws.Range("H61").Value = 1859.0625
$ws.Range("I61").Value = 2340
$ws.Range("J61").Value = 1790.3572
$ws.Range("K61").Value = 2340
$ws.Range("L61").Value = 1790.3572
$ws.Range("M61").Value = -2138
$ws.Range("N61").Value = -2194.3572
$ws.Range("H113").Value = 1859.0625
$ws.Range("I113").Value = 2340
$ws.Range("J113").Value = 1790.3572
$ws.Range("K113").Value = 2340
$ws.Range("L113").Value = 1790.3572
$ws.Range("M113").Value = -170
$ws.Range("N113").Value = -6130.3572
$ws.Range("H136").Value = 2704.049
$ws.Range("I136").Value = 1002.24
$ws.Range("J136").Value = 10439.546
$ws.Range("K136").Value = 3006.72
$ws.Range("L136").Value = 31318.638
$ws.Range("M136").Value = -456.7200000000003
$ws.Range("N136").Value = -36418.638
